# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect the latest generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 710
$wsExhibit.Range("F5").Value = 2439
$wsExhibit.Range("F6").Value = 53
$wsExhibit.Range("F7").Value = 3563
$wsExhibit.Range("F8").Value = 466
$wsExhibit.Range("F9").Value = 905

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 710
$wsAll.Range("F6").Value = 2439
$wsAll.Range("F7").Value = 53
$wsAll.Range("F8").Value = 3563
$wsAll.Range("F9").Value = 466
$wsAll.Range("F10").Value = 905
